$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (data index 0-23)
# Columns: B, D, E, F, G, L (col indices 2, 4, 5, 6, 7, 12)

$colB = @(20.93624178807288, 20.54987272643285, 20.31469454537249, 20.21949869737812, 20.20373387657069, 20.31340793968598, 20.80266998673138, 21.77299625187807, 22.48546241217728, 22.80806585317183, 22.92989958922019, 22.90367689685235, 22.81809641369756, 22.76562976309735, 22.46433949348061, 22.27904076681869, 22.17232699111524, 22.13617592786397, 22.29878092851833, 22.8432432841767, 23.1971138368016, 23.00846235354179, 22.28985696296986, 21.51007494501343)
$colD = @(8.007321480904176, 7.958596390032397, 7.933830317638796, 7.925020047863359, 7.923634093449548, 7.933706327979796, 7.989444977257598, 8.140044548018139, 8.27620527094952, 8.343692194333924, 8.370040837897786, 8.364331032101029, 8.345844061523714, 8.334623314571008, 8.271906111329043, 8.234849294359353, 8.214057757708922, 8.207107950450974, 8.238739998043524, 8.351252668612895, 8.429403139700414, 8.387272617713453, 8.236979412941087, 8.094834061945555)
$colE = @(26.42599522793639, 24.86685746047598, 23.85661798970947, 23.43174259900401, 23.36039945639675, 23.85094122075382, 25.89944597340774, 29.49330409754084, 31.87306352521066, 32.89893113702733, 33.27925190359166, 33.1977055074591, 32.93038371289065, 32.76557982846833, 31.80487949450124, 31.20099555406181, 30.84832005691125, 30.72799368761882, 31.26583247227232, 33.00912381772223, 34.10095119282564, 33.52256702371124, 31.23653681579315, 28.56673599358854)
$colF = @(45.73867726674013, 44.67071753454178, 44.03141973765547, 43.77540837213416, 43.73318089601745, 44.02794833885007, 45.367261604772, 48.10728683143748, 50.1673711829193, 51.11022098951422, 51.4677174654156, 51.39070928074593, 51.13962521951782, 50.98587849293978, 50.10583725103375, 49.56715899844269, 49.2578831096281, 49.1532747046919, 49.62444718915462, 51.21336499108013, 52.25434211063816, 51.69863548678651, 49.59854587034846, 47.35619849320319)
$colG = @(3.654038404568322, 3.663261123953052, 3.669183571179548, 3.671662807201172, 3.672078470489768, 3.669216739983044, 3.657164807446974, 3.635568377312479, 3.620911034756401, 3.614498505104904, 3.612106356547403, 3.612619950302316, 3.614300980587685, 3.615335349615454, 3.621335196059589, 3.625080870613514, 3.627259335183057, 3.628001073490966, 3.62467965189937, 3.613806244973591, 3.606910214154971, 3.61057169348294, 3.624860964886804, 3.641195992914038)
$colL = @(13.38062912091392, 13.03367191781584, 12.81956719320399, 12.73217721965457, 12.71766156432649, 12.81838901624028, 13.26129625272413, 14.11598174848211, 14.72880994551854, 15.00306743614303, 15.10617814692831, 15.08400594441004, 15.01156585173559, 14.96709456657238, 14.71078732514815, 14.55232324396755, 14.46075846231833, 14.4296870085561, 14.5692362405801, 15.03286418399541, 15.33149354150994, 15.17253910342913, 14.56159130554833, 13.886967354046)

$columns = @{
    2 = $colB
    4 = $colD
    5 = $colE
    6 = $colF
    7 = $colG
    12 = $colL
}

foreach ($colIndex in $columns.Keys) {
    $values = $columns[$colIndex]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $rowNum = $i + 2
        $ws.Cells.Item($rowNum, $colIndex).Value = $values[$i]
    }
}
